# Onsen UI quiz workbook — "Add files via upload" edit
# Applies the content changes: a few answer corrections/re-highlights,
# two new trailing quiz rows, refreshed selection, and print setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Copy-CellFormat($srcAddress, $dstAddress) {
    $ws.Range($srcAddress).Copy() | Out-Null
    $ws.Range($dstAddress).PasteSpecial($xlPasteFormats) | Out-Null
}

# --- Row 19: "Which of the following is the correct way of specifying an
#     event handler?" -> answer text unchanged ("onClick"); re-highlighted
#     red/bordered (style matches C40/C43/.../C50).
Copy-CellFormat "C40" "C19"

# --- Row 22: "_____ helps in popping the top page of the navigator."
#     answer corrected from "Popbutton" to "Pushbutton".
$ws.Range("C22").Value = "Pushbutton"

# --- Row 25: "The <Carosal> component can render content horizontally and
#     vertically" -> answer unchanged (TRUE); re-highlighted green/bordered
#     (style matches C2/C4/...).
Copy-CellFormat "C2" "C25"

# --- Row 35: "The <ToolBar> component will render a toolbar at the _____ of
#     the page:" answer corrected from "Top/left" to "Top" and
#     re-highlighted green/bordered.
$ws.Range("C35").Value = "Top"
Copy-CellFormat "C2" "C35"

# --- Row 53: "The prop that decides whether a dialog can be cancelled :"
#     answer corrected from "oncancel" to "isCancelable" and re-highlighted
#     green/no-border (style matches C48).
$ws.Range("C53").Value = "isCancelable"
Copy-CellFormat "C48" "C53"

# --- Two new quiz rows appended at the bottom of the sheet.
$ws.Range("B55").Value = "The ons.platform,select() should be called after the app initializes completely"
$ws.Range("C55").Value = $false
Copy-CellFormat "C48" "C55"

$ws.Range("B56").Value = "The following are built-in animation interfaces, except"
$ws.Range("C56").Value = "SwiperAnimator"
Copy-CellFormat "C48" "C56"

# --- Refresh the view: scroll position + active selection.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C22").Select() | Out-Null

# --- Print setup (paper size / orientation) picked up by the saved file.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
